$p = $ppt.ActivePresentation

# Insert the new 4th slide ("Git add") using the same "Title and Content" layout
# that the existing slides (e.g. slide 3, "Git status") already use.
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title placeholder ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = 'Git'
$title.LanguageID = "es-GT"
$title.ParagraphFormat.Alignment = 2  # ppAlignCenter -> <a:pPr algn="ctr"/>
$titleRun = $title
$titleRun = $titleRun.InsertAfter(' ')
$titleRun.LanguageID = "es-GT"
$titleRun = $titleRun.InsertAfter('add')
$titleRun.LanguageID = "es-GT"

# --- Content placeholder ---
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = 'El comando '
$body.LanguageID = "es-GT"
$bodyRun = $body
$bodyRun = $bodyRun.InsertAfter('git')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter(' ')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter('add')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter(' añade contenido del directorio de trabajo al área de ensayo (')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter('staging')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter(' ')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter('area')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter(' o ''')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter('index')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter(''') para la próxima confirmación. Cuando se ejecuta el comando ')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter('git')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter(' ')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter('commit')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter(', éste, de forma predeterminada, sólo mira en esta área de ensayo, por lo que ')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter('git')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter(' ')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter('add')
$bodyRun.LanguageID = "es-GT"
$bodyRun = $bodyRun.InsertAfter(' se utiliza para fabricar exactamente lo que te gustaría fuese tu próxima instantánea a confirmar.')
$bodyRun.LanguageID = "es-GT"

